$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 191.66667
$ws.Range("I12").Value = 210
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 210
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -40
$ws.Range("N12").Value = -440
$ws.Range("H33").Value = 263.7
$ws.Range("J33").Value = 549.3333
$ws.Range("L33").Value = 549.3333
$ws.Range("N33").Value = -1007.3333
$ws.Range("H37").Value = 950
$ws.Range("I37").Value = 900
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 2700
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = -2574
$ws.Range("N37").Value = -3252
$ws.Range("H74").Value = 2900
$ws.Range("I74").Value = 2800
$ws.Range("K74").Value = 2800
$ws.Range("M74").Value = -1864
$ws.Range("H77").Value = 2900
$ws.Range("I77").Value = 2800
$ws.Range("K77").Value = 14000
$ws.Range("M77").Value = -9320
$ws.Range("H100").Value = 1771.1818
$ws.Range("I100").Value = 1440
$ws.Range("J100").Value = 2000.4615
$ws.Range("K100").Value = 1440
$ws.Range("L100").Value = 2000.4615
$ws.Range("M100").Value = -899
$ws.Range("N100").Value = -3082.4615
$ws.Range("H113").Value = 3501.111
$ws.Range("I113").Value = 3142
$ws.Range("J113").Value = 3950
$ws.Range("K113").Value = 3142
$ws.Range("L113").Value = 3950
$ws.Range("M113").Value = 112
$ws.Range("N113").Value = -10458
$ws.Range("H132").Value = 7097814.5
$ws.Range("I132").Value = 10758095
$ws.Range("K132").Value = 32274285
$ws.Range("M132").Value = -32271755
$ws.Range("H135").Value = 657.44446
$ws.Range("I135").Value = 424.875
$ws.Range("K135").Value = 3823.875
$ws.Range("M135").Value = -1288.875
$ws.Range("H138").Value = 1395.0404
$ws.Range("I138").Value = 884.8570999999999
$ws.Range("J138").Value = 1532.3975
$ws.Range("K138").Value = 2654.5713
$ws.Range("L138").Value = 4597.1925
$ws.Range("M138").Value = 2485.4287
$ws.Range("N138").Value = -14877.1925
$ws.Range("H141").Value = 532.2692
$ws.Range("I141").Value = 551.2917
$ws.Range("J141").Value = 304
$ws.Range("K141").Value = 1653.8751
$ws.Range("L141").Value = 912
$ws.Range("M141").Value = 3526.1249
$ws.Range("N141").Value = -11272

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2505.2468
$ws.Range("I32").Value = 2272.0532
$ws.Range("K32").Value = 2272.0532
$ws.Range("M32").Value = -1985.0532
$ws.Range("H61").Value = 1397.7106
$ws.Range("I61").Value = 1296.1666
$ws.Range("J61").Value = 1778.5
$ws.Range("K61").Value = 1296.1666
$ws.Range("L61").Value = 1778.5
$ws.Range("M61").Value = -1084.1666
$ws.Range("N61").Value = -2202.5
$ws.Range("H102").Value = 83334430
$ws.Range("I102").Value = 83334430
$ws.Range("K102").Value = 83334430
$ws.Range("M102").Value = -83332808
$ws.Range("H110").Value = 1580.7391
$ws.Range("I110").Value = 1188.9286
$ws.Range("J110").Value = 2190.2222
$ws.Range("K110").Value = 1188.9286
$ws.Range("L110").Value = 2190.2222
$ws.Range("M110").Value = 856.0714
$ws.Range("N110").Value = -6280.2222
$ws.Range("H132").Value = 1542.4694
$ws.Range("I132").Value = 1321.0883
$ws.Range("J132").Value = 2044.2667
$ws.Range("K132").Value = 3963.2649
$ws.Range("L132").Value = 6132.800099999999
$ws.Range("M132").Value = -1433.2649
$ws.Range("N132").Value = -11192.8001
$ws.Range("H136").Value = 1397.7106
$ws.Range("I136").Value = 1296.1666
$ws.Range("J136").Value = 1778.5
$ws.Range("K136").Value = 3888.4998
$ws.Range("L136").Value = 5335.5
$ws.Range("M136").Value = -1338.4998
$ws.Range("N136").Value = -10435.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2688.0588
$ws.Range("I86").Value = 3520.6
$ws.Range("J86").Value = 1498.7142
$ws.Range("K86").Value = 3520.6
$ws.Range("L86").Value = 1498.7142
$ws.Range("M86").Value = -2397.6
$ws.Range("N86").Value = -3744.7142
$ws.Range("H89").Value = 2688.0588
$ws.Range("I89").Value = 3520.6
$ws.Range("J89").Value = 1498.7142
$ws.Range("K89").Value = 17603
$ws.Range("L89").Value = 7493.571
$ws.Range("M89").Value = -11987
$ws.Range("N89").Value = -18725.571
$ws.Range("H94").Value = 16667644
$ws.Range("I94").Value = 17857976
$ws.Range("J94").Value = 3010
$ws.Range("K94").Value = 17857976
$ws.Range("L94").Value = 3010
$ws.Range("M94").Value = -17857525
$ws.Range("N94").Value = -3912
$ws.Range("H99").Value = 76924240
$ws.Range("I99").Value = 90910216
$ws.Range("J99").Value = 1350
$ws.Range("K99").Value = 90910216
$ws.Range("L99").Value = 1350
$ws.Range("M99").Value = -90908718
$ws.Range("N99").Value = -4346
$ws.Range("H105").Value = 55557132
$ws.Range("I105").Value = 62501624
$ws.Range("J105").Value = 1205.5
$ws.Range("K105").Value = 62501624
$ws.Range("L105").Value = 1205.5
$ws.Range("M105").Value = -62499877
$ws.Range("N105").Value = -4699.5
$ws.Range("H132").Value = 1712666.5
$ws.Range("J132").Value = 1712666.5
$ws.Range("L132").Value = 1712666.5
$ws.Range("N132").Value = -1722786.5
$ws.Range("H134").Value = 3367.1
$ws.Range("I134").Value = 799.4048
$ws.Range("J134").Value = 16847.5
$ws.Range("K134").Value = 2398.2144
$ws.Range("L134").Value = 50542.5
$ws.Range("M134").Value = 136.7856000000002
$ws.Range("N134").Value = -55612.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 28573486
$ws.Range("I62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("H65").Value = 28573486
$ws.Range("I65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("H134").Value = 1080.7188
$ws.Range("I134").Value = 1091.6923
$ws.Range("J134").Value = 1033.1666
$ws.Range("K134").Value = 3275.0769
$ws.Range("L134").Value = 3099.4998
$ws.Range("M134").Value = -740.0769
$ws.Range("N134").Value = -8169.4998

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 535004.8
$ws.Range("I4").Value = 141.66667
$ws.Range("J4").Value = 1118491.9
$ws.Range("K4").Value = 425.00001
$ws.Range("L4").Value = 3355475.7
$ws.Range("M4").Value = -313.00001
$ws.Range("N4").Value = -3355699.7
$ws.Range("H56").Value = 5710.52
$ws.Range("I56").Value = 5710.52
$ws.Range("K56").Value = 5710.52
$ws.Range("M56").Value = -5180.52
$ws.Range("H64").Value = 3856.8857
$ws.Range("I64").Value = 4012
$ws.Range("J64").Value = 3852.3235
$ws.Range("K64").Value = 12036
$ws.Range("L64").Value = 11556.9705
$ws.Range("M64").Value = -11766
$ws.Range("N64").Value = -12096.9705
$ws.Range("H67").Value = 3856.8857
$ws.Range("I67").Value = 4012
$ws.Range("J67").Value = 3852.3235
$ws.Range("K67").Value = 12036
$ws.Range("L67").Value = 11556.9705
$ws.Range("M67").Value = -11100
$ws.Range("N67").Value = -13428.9705
$ws.Range("H140").Value = 22715.055
$ws.Range("I140").Value = 52384.5
$ws.Range("J140").Value = 2935.4243
$ws.Range("K140").Value = 157153.5
$ws.Range("L140").Value = 8806.2729
$ws.Range("M140").Value = -151973.5
$ws.Range("N140").Value = -19166.2729

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 679.875
$ws.Range("I97").Value = 705.5714
$ws.Range("K97").Value = 705.5714
$ws.Range("M97").Value = -209.5714
$ws.Range("H132").Value = 1993.0667
$ws.Range("I132").Value = 1487.1
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 4461.299999999999
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -1931.299999999999
$ws.Range("N132").Value = -14075

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1180.5385
$ws.Range("I136").Value = 901.4211
$ws.Range("J136").Value = 1938.1428
$ws.Range("K136").Value = 2704.2633
$ws.Range("L136").Value = 5814.428400000001
$ws.Range("M136").Value = -154.2633000000001
$ws.Range("N136").Value = -10914.4284

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 23814750
$ws.Range("I62").Value = 29417342
$ws.Range("J62").Value = 3737.5
$ws.Range("K62").Value = 29417342
$ws.Range("L62").Value = 3737.5
$ws.Range("M62").Value = -29416718
$ws.Range("N62").Value = -4985.5
$ws.Range("H65").Value = 23814750
$ws.Range("I65").Value = 29417342
$ws.Range("J65").Value = 3737.5
$ws.Range("K65").Value = 147086710
$ws.Range("L65").Value = 18687.5
$ws.Range("M65").Value = -147083590
$ws.Range("N65").Value = -24927.5
$ws.Range("H102").Value = 14000
$ws.Range("J102").Value = 14000
$ws.Range("L102").Value = 14000
$ws.Range("N102").Value = -20490
$ws.Range("H132").Value = 2066.8845
$ws.Range("I132").Value = 2071.628
$ws.Range("J132").Value = 2044.2222
$ws.Range("K132").Value = 6214.884
$ws.Range("L132").Value = 6132.6666
$ws.Range("M132").Value = -3684.884
$ws.Range("N132").Value = -11192.6666
